$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CTP")

# Insert two new rows above row 13 (one at a time, so each is a clean
# single-row insert that shifts everything below down by one).
$ws.Range("A13:E13").EntireRow.Insert()
$ws.Range("A13:E13").EntireRow.Insert()

# The freshly inserted rows pick up "no border" formatting; restore the
# left/right boundary-box formatting (style of column A / E) by copying
# it down from the row that still carries it (now row 16, formerly 14).
$ws.Range("A16").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E13:E14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new rows with the new request/response pair of entries.
$ws.Range("C13").Value = "ReqQryTradingAccount"
$ws.Range("D13").Value = "OnRspQryTradingAccount"
$ws.Range("C14").Value = "ReqQryInvestorPosition"
$ws.Range("D14").Value = "onRspQryInvestorPosition"

# Match the author's final selection in the saved workbook.
$ws.Range("D14").Select()
